# Word COM-interop script implementing the commit "Updated for Temp folder".
#
# Semantic change: insert a new bullet list item ("Temp folder for saving
# the uploaded files is given as D:/temp/ by default in
# application.properties. Make sure it is present, or change as per
# need.") right after the "Make sure Maven and Java is set in the
# current path" bullet and before the "Run the command mvn
# spring-boot:run" bullet. A handful of small proofing/cosmetic tweaks
# ride along with it (the grammar/spell-check re-scan that Word performs
# after an edit moves the hidden "_GoBack" bookmark to the new insertion
# point, merges a couple of runs that used to be split around
# now-resolved grammar flags, and the picture run picks up an explicit
# language tag).

$d = $word.ActiveDocument

# --- 1. Insert the new "Temp folder" list paragraph ------------------------
# Paragraph 2 is "Make sure Maven and Java is set in the current path".
$mavenPara = $d.Paragraphs(2)
$newPara = $mavenPara.Range.InsertParagraphAfter()

# The freshly created paragraph (Paragraphs(3)) already carries the
# ListParagraph style / numPr, because it was split off the Maven
# paragraph, matching the bullet list it belongs to.
$tempPara = $d.Paragraphs(3)
$tempPara.Range.Text = "Temp folder for saving the uploaded files is given as D:/temp/ by default in application.properties. Make sure it is present, or change as per need."

# --- 2. Move the "_GoBack" bookmark to the new insertion point -------------
# Word stamps "_GoBack" at the location of the most recent edit; here that
# is right after the user finished typing "application.propertie" (just
# before the final "s" of "application.properties").
$bookmarkRange = $d.Content.Duplicate
$bookmarkRange.Find.Execute("application.propertie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 3. Small proofing clean-ups that ride along with the re-scan ----------

# "PostMan,etc" + "." + "," -> "PostMan,etc" + ".,"  (the trailing period
# and comma, previously split apart by a now-stale grammar flag, collapse
# back into one run).
$d.Content.Find.Execute(".,", $true, $false, $false, $false, $false, $true, 1, $false, ".,", 2)

# "has to" / " be checked for each line equality" -> one run, the
# now-resolved grammar flag drops out.
$d.Content.Find.Execute("has to be checked for each line equality", $true, $false, $false, $false, $false, $true, 1, $false, "has to be checked for each line equality", 2)

# "Return type of the api is json " : split "json" back into its own run
# (re-flagged by the spell checker).
$jsonRange = $d.Content.Duplicate
$jsonRange.Find.Execute("json", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$jsonRange.Bold = 1
$jsonRange.Bold = 0

# The screenshot picture run picks up an explicit proofing language (the
# author's Word install is set to English (India)).
$picturePara = $d.Paragraphs(6)
$pictureRange = $picturePara.Range.InlineShapes(1).Range
$pictureRange.LanguageID = "en-IN"
$pictureRange.LanguageIDFarEast = "en-IN"

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
foreach ($p in $d.Paragraphs) {
    Write-Output $p.Range.Text
}
